$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range("A32").Value = 111357523
$ws.Range("B32").Value = 78081
$ws.Range("E32").Value = 229821
$ws.Range("F32").Value = "Vedflamlav"
$ws.Range("G32").Value = "Ramboldia elabens"
$ws.Range("H32").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q32").Value = 367297.1312965819
$ws.Range("R32").Value = 6877354.545063579
$ws.Range("Z32").Value = "00:00"
$ws.Range("AB32").Value = "00:00"
$ws.Range("AC32").Value = "växer på silverstubbe  med brandljud"
$ws.Range("AI32").Value = "Naturskog."
$ws.Range("AJ32").Value = "tall"
$ws.Range("AK32").Value = "Pinus sylvestris"
$ws.Range("AO32").Value = "Pinus sylvestris"

# Row 33
$ws.Range("A33").Value = 111358214
$ws.Range("B33").Value = 78579
$ws.Range("E33").Value = 2081
$ws.Range("F33").Value = "Skrovellav"
$ws.Range("G33").Value = "Lobaria scrobiculata"
$ws.Range("H33").Value = "(Scop.) DC."
$ws.Range("Q33").Value = 367225.8747160842
$ws.Range("R33").Value = 6877314.542789092
$ws.Range("Z33").Value = "13:30"
$ws.Range("AB33").Value = "13:30"
$ws.Range("AC33").Value = "Växer på gamal sälg"
$ws.Range("AI33").Value = "Fuktig grannaturskog"
$ws.Range("AJ33").Value = "sälg"
$ws.Range("AK33").Value = "Salix caprea"
$ws.Range("AO33").Value = "Salix caprea"

# Row 34
$ws.Range("A34").Value = 111357955
$ws.Range("I34").Value = "30"
$ws.Range("J34").Value = "plantor/tuvor"
$ws.Range("Q34").Value = 367250.4893208001
$ws.Range("R34").Value = 6877317.812129297
$ws.Range("Z34").Value = "13:00"
$ws.Range("AB34").Value = "13:00"
$ws.Range("AC34").Value = ""

# Row 35
$ws.Range("A35").Value = 111355912
$ws.Range("B35").Value = 73696
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 6440
$ws.Range("F35").Value = "Vitgrynig nållav"
$ws.Range("G35").Value = "Chaenotheca subroscida"
$ws.Range("H35").Value = "(Eitner) Zahlbr."
$ws.Range("L35").Value = ""
$ws.Range("Q35").Value = 367365.3834173826
$ws.Range("R35").Value = 6877426.168084201
$ws.Range("Z35").Value = "12:00"
$ws.Range("AB35").Value = "12:00"
$ws.Range("AC35").Value = "Växer på gammal gran som är snitslad av Mellanskog."
$ws.Range("AI35").Value = "Grannaturskog"
$ws.Range("AJ35").Value = "gran"
$ws.Range("AK35").Value = "Picea abies"
$ws.Range("AO35").Value = "Picea abies"

# Row 36
$ws.Range("A36").Value = 111357675
$ws.Range("I36").Value = ""
$ws.Range("J36").Value = ""
$ws.Range("Q36").Value = 367367.2640777439
$ws.Range("R36").Value = 6877426.094615285
$ws.Range("Z36").Value = "12:00"
$ws.Range("AB36").Value = "12:00"
$ws.Range("AC36").Value = "Många plantor av spindelblosmter. Finns utspridd i hela skogen."

# Row 37
$ws.Range("A37").Value = 111357840
$ws.Range("B37").Value = 77515
$ws.Range("D37").Value = "NT"
$ws.Range("E37").Value = 6425
$ws.Range("F37").Value = "Garnlav"
$ws.Range("G37").Value = "Alectoria sarmentosa"
$ws.Range("H37").Value = "(Ach.) Ach."
$ws.Range("L37").Value = ""
$ws.Range("Q37").Value = 367293.7760285549
$ws.Range("R37").Value = 6877401.225212204
$ws.Range("Z37").Value = "12:30"
$ws.Range("AB37").Value = "12:30"
$ws.Range("AC37").Value = "Skogen är hänglavsrik med långa bålar."

# Row 38
$ws.Range("A38").Value = 111358801
$ws.Range("Q38").Value = 367379.1282350773
$ws.Range("R38").Value = 6877428.452265346
$ws.Range("Z38").Value = "00:00"
$ws.Range("AB38").Value = "00:00"
$ws.Range("AC38").Value = "10 -tals plantor"
$ws.Range("AI38").Value = "Fuktig grannaturskog"

# Row 39
$ws.Range("A39").Value = 111358171
$ws.Range("B39").Value = 96368
$ws.Range("D39").Value = "LC"
$ws.Range("E39").Value = 221952
$ws.Range("F39").Value = "Spindelblomster"
$ws.Range("G39").Value = "Neottia cordata"
$ws.Range("H39").Value = "(L.) Rich."
$ws.Range("L39").Value = ""
$ws.Range("Q39").Value = 367256.9250051867
$ws.Range("R39").Value = 6877313.798904967
$ws.Range("Z39").Value = "13:00"
$ws.Range("AB39").Value = "13:00"
$ws.Range("AC39").Value = "Flertalet plantor"
$ws.Range("AJ39").Value = ""
$ws.Range("AK39").Value = ""
$ws.Range("AO39").Value = ""

# Row 40
$ws.Range("A40").Value = 111357794
$ws.Range("B40").Value = 96368
$ws.Range("D40").Value = "LC"
$ws.Range("E40").Value = 221952
$ws.Range("F40").Value = "Spindelblomster"
$ws.Range("G40").Value = "Neottia cordata"
$ws.Range("H40").Value = "(L.) Rich."
$ws.Range("L40").Value = ""
$ws.Range("Q40").Value = 367349.0144895007
$ws.Range("R40").Value = 6877453.137823995
$ws.Range("Z40").Value = "11:30"
$ws.Range("AB40").Value = "11:30"
$ws.Range("AC40").Value = "Flertalet plantor inom några kvm i grannaturskog"

# Row 41
$ws.Range("A41").Value = 111358288
$ws.Range("B41").Value = 95532
$ws.Range("D41").Value = "LC"
$ws.Range("E41").Value = 221945
$ws.Range("F41").Value = "Revlummer"
$ws.Range("G41").Value = "Lycopodium annotinum"
$ws.Range("H41").Value = "L."
$ws.Range("L41").Value = ""
$ws.Range("Q41").Value = 367244.7553745224
$ws.Range("R41").Value = 6877315.685314978
$ws.Range("AC41").Value = "Lummerväxter växer i stor mängd i skogsområdet"
$ws.Range("AI41").Value = "Grannaturskog"
$ws.Range("AJ41").Value = ""
$ws.Range("AK41").Value = ""
$ws.Range("AO41").Value = ""

# Row 42
$ws.Range("A42").Value = 111358524
$ws.Range("B42").Value = 78579
$ws.Range("E42").Value = 2081
$ws.Range("F42").Value = "Skrovellav"
$ws.Range("G42").Value = "Lobaria scrobiculata"
$ws.Range("H42").Value = "(Scop.) DC."
$ws.Range("Q42").Value = 367239.6455721884
$ws.Range("R42").Value = 6877329.520917758
$ws.Range("Z42").Value = "13:30"
$ws.Range("AB42").Value = "13:30"
$ws.Range("AC42").Value = "Växer på gammal sälg i fuktig grannaturskog"
$ws.Range("AI42").Value = "Fuktig grannaturskog"
$ws.Range("AJ42").Value = "sälg"
$ws.Range("AK42").Value = "Salix caprea"
$ws.Range("AO42").Value = "Salix caprea"

# Row 43
$ws.Range("A43").Value = 111358759
$ws.Range("B43").Value = 78579
$ws.Range("D43").Value = "NT"
$ws.Range("E43").Value = 2081
$ws.Range("F43").Value = "Skrovellav"
$ws.Range("G43").Value = "Lobaria scrobiculata"
$ws.Range("H43").Value = "(Scop.) DC."
$ws.Range("L43").Value = ""
$ws.Range("Q43").Value = 367361.6488811045
$ws.Range("R43").Value = 6877475.21281712
$ws.Range("Z43").Value = "11:00"
$ws.Range("AB43").Value = "11:00"
$ws.Range("AC43").Value = "Växer på grov sälg"
$ws.Range("AJ43").Value = "sälg"
$ws.Range("AK43").Value = "Salix caprea"
$ws.Range("AO43").Value = "Salix caprea"

# Row 44
$ws.Range("A44").Value = 111366550
$ws.Range("M44").Value = ""
$ws.Range("Q44").Value = 367370.4653471896
$ws.Range("R44").Value = 6877327.229400039
$ws.Range("Z44").Value = "14:00"
$ws.Range("AB44").Value = "14:00"
$ws.Range("AC44").Value = "Hoppade i blåbärsriset ner i en källhåla i skogen"
$ws.Range("AI44").Value = "Grannaturskog"

# Row 46
$ws.Range("A46").Value = 111366200
$ws.Range("M46").Value = "födosökande"
$ws.Range("Q46").Value = 367379.213297962
$ws.Range("R46").Value = 6877418.575064567
$ws.Range("Z46").Value = "13:30"
$ws.Range("AB46").Value = "13:30"
$ws.Range("AC46").Value = "Stor vanlig groda i bäckfåra. Flertalet grodor observerades i skogsområdet"
$ws.Range("AI46").Value = "Källpåverkad grannaturskog"
